# "Devdutt Padikkal" innings log: the per-innings batting stats (runs/balls/
# fours/sixes in columns C:F) are reordered into a new row sequence - commit
# message: "updated activity till excel form". Every value already exists
# somewhere in the sheet; this just moves each innings to its new row.
# playerName/teamName (A:B), the header row, and row 12 (unchanged innings)
# are left untouched, and a cell is only rewritten when its value actually
# changes (its runs/balls/fours/sixes can coincidentally repeat row to row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: target row -> new "runs","balls","fours","sixes" text values for C:F
$newRowValues = @{
    2 = @("1", "6", "0", "0")   # was row 3
    3 = @("50", "41", "5", "0")   # was row 5
    4 = @("22", "21", "2", "1")   # was row 9
    5 = @("5", "8", "1", "0")   # was row 4
    6 = @("74", "45", "12", "1")   # was row 2
    7 = @("32", "23", "4", "1")   # was row 14
    8 = @("1", "2", "0", "0")   # was row 13
    9 = @("18", "12", "1", "1")   # was row 10
    10 = @("63", "45", "6", "1")   # was row 7
    11 = @("33", "34", "2", "1")   # was row 6
    13 = @("54", "40", "5", "2")   # was row 11
    14 = @("25", "17", "3", "0")   # was row 8
    15 = @("56", "42", "8", "0")   # was row 16
    16 = @("4", "6", "0", "0")   # was row 15
}

$cols = @("C", "D", "E", "F")
foreach ($r in $newRowValues.Keys) {
    $vals = $newRowValues[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cellRef = "$($cols[$i])$r"
        $cell = $ws.Range($cellRef)
        $newVal = $vals[$i]
        if ($cell.Text -ne $newVal) {
            # Keep these numeric-looking values stored as text, matching the
            # sheet-wide "number stored as text" state already on this range.
            $cell.NumberFormat = "@"
            $cell.Value = $newVal
        }
    }
}
